$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C. This shifts the existing "Email ID",
# "Year" and "Due for receipt" columns one position to the right.
$ws.Range("C1").EntireColumn.Insert()

# Populate the new column's header with the "Employee Status" label and
# give it the same centered header formatting used by the other headers
# (e.g. "Sr. No.", "Name", "Email ID", "Due for receipt").
$ws.Range("C1").Value = "Employee Status"
$ws.Range("C1").HorizontalAlignment = -4108

# Match the resulting selection state of the workbook.
$ws.Range("C2").Select()
